$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = '@'
$r.Value = '56.306.65'
$r.ClearFormats()
$ws.Range('E2').Value = '  +3.51%  '
$r = $ws.Range('D3')
$r.NumberFormat = '@'
$r.Value = '2.966.63'
$r.ClearFormats()
$ws.Range('E3').Value = '  +2.40%  '
$ws.Range('E4').Value = '  +0.06%  '
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '500.30'
$r.ClearFormats()
$ws.Range('E5').Value = '  +4.83%  '
$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '133.97'
$r.ClearFormats()
$ws.Range('E6').Value = '  +5.10%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +5.61%  '
$r = $ws.Range('D9')
$r.NumberFormat = '@'
$r.Value = '7.34'
$r.ClearFormats()
$ws.Range('E9').Value = '  +9.78%  '
$ws.Range('E10').Value = '  +8.29%  '
$ws.Range('E11').Value = '  +3.89%  '
$ws.Range('E12').Value = '  +2.95%  '
$r = $ws.Range('D13')
$r.NumberFormat = '@'
$r.Value = '3.475.55'
$r.ClearFormats()
$ws.Range('E13').Value = '  +2.95%  '
$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '25.16'
$r.ClearFormats()
$ws.Range('E14').Value = '  +9.41%  '
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '56.291.39'
$r.ClearFormats()
$ws.Range('E15').Value = '  +3.64%  '
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '0.0000150'
$r.ClearFormats()
$ws.Range('E16').Value = '  +10.02%  '
$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '2.966.70'
$r.ClearFormats()
$ws.Range('E17').Value = '  +3.03%  '
$ws.Range('E18').Value = '  +7.76%  '
$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '12.25'
$r.ClearFormats()
$ws.Range('E19').Value = '  +4.78%  '
$ws.Range('E20').Value = '  +7.36%  '
$r = $ws.Range('D21')
$r.NumberFormat = '@'
$r.Value = '320.02'
$r.ClearFormats()
$ws.Range('E21').Value = '  +2.65%  '
$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '1.00'
$r.ClearFormats()
$ws.Range('E22').Value = '  -0.27%  '
$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '0.466'
$r.ClearFormats()
$ws.Range('E23').Value = '  +3.43%  '
$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '61.76'
$r.ClearFormats()
$ws.Range('E24').Value = '  +3.18%  '
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '0.997'
$r.ClearFormats()
$ws.Range('E25').Value = '  -0.38%  '
$ws.Range('E26').Value = '  +4.29%  '
$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '0.0₃0879'
$r.ClearFormats()
$ws.Range('E27').Value = '  +5.14%  '
$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '6.36'
$r.ClearFormats()
$ws.Range('E28').Value = '  -0.32%  '
$ws.Range('E29').Value = '  +7.44%  '
$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '1.17'
$r.ClearFormats()
$ws.Range('E30').Value = '  +0.87%  '
$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '1.73'
$r.ClearFormats()
$ws.Range('E31').Value = '  +6.45%  '
$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '20.19'
$r.ClearFormats()
$ws.Range('E32').Value = '  +4.76%  '
$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '158.28'
$r.ClearFormats()
$ws.Range('E33').Value = '  +13.00%  '
$r = $ws.Range('D34')
$r.NumberFormat = '@'
$r.Value = '4.40'
$r.ClearFormats()
$ws.Range('E34').Value = '  +2.52%  '
$ws.Range('E35').Value = '  +2.11%  '
$r = $ws.Range('D36')
$r.NumberFormat = '@'
$r.Value = '5.50'
$r.ClearFormats()
$ws.Range('E36').Value = '  +0.09%  '
$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '0.0668'
$r.ClearFormats()
$ws.Range('E37').Value = '  +6.70%  '
$ws.Range('E38').Value = '  -0.78%  '
$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '3.001.50'
$r.ClearFormats()
$ws.Range('E39').Value = '  +3.19%  '
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '0.999'
$r.ClearFormats()
$ws.Range('E40').Value = '  +0.12%  '
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '36.07'
$r.ClearFormats()
$ws.Range('E41').Value = '  +2.02%  '
$ws.Range('E42').Value = '  +5.98%  '
$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '2.226.93'
$r.ClearFormats()
$ws.Range('E43').Value = '  +7.89%  '
$ws.Range('E44').Value = '  +3.59%  '
$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '0.974'
$r.ClearFormats()
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('E46').Value = '  +1.88%  '
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '1.91'
$r.ClearFormats()
$ws.Range('E47').Value = '  +15.56%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '0.0233'
$r.ClearFormats()
$ws.Range('E48').Value = '  +8.58%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$r = $ws.Range('D49')
$r.NumberFormat = '@'
$r.Value = '5.72'
$r.ClearFormats()
$ws.Range('E49').Value = '  +6.55%  '
$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '18.76'
$r.ClearFormats()
$ws.Range('E50').Value = '  +2.95%  '
$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '0.0861'
$r.ClearFormats()
$ws.Range('E51').Value = '  +7.09%  '
